# Add a "timezone" / "Asia/Qatar" row to the tournament info table on the
# Tournament sheet (new row 5; existing rows 5-12 shift down to 6-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tournament")
$ws.Activate()

# Insert a blank worksheet row above the old row 5 ("venue.1"), shifting the
# remaining table rows down.
$ws.Rows.Item(5).Insert()

# Grow the "tournament" table definition to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I13"))

# Populate the new row: key = "timezone", en = "Asia/Qatar" (other language
# columns are intentionally left blank, same as the source edit).
$ws.Range("A5").Value = "timezone"
$ws.Range("B5").Value = "Asia/Qatar"

# Leave the selection on the newly-added row, matching the author's saved
# view state.
$ws.Rows.Item(5).Select()
